# Applies the latest cryptos-list scrape: updates Price (D) / Volume(1h) (E)
# columns for each coin row, plus a Stellar/VeChain row-order swap (rows 47-48).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: the target cell, its new text, and whether the text must be
# force-written as Text (some Price values, e.g. "598.40" or "0.999", would
# otherwise be auto-coerced into numbers by Excel and lose their formatting).
$updates = @(
    @{ Cell = 'D2'; Value = '65.795.65'; ForceText = $false },
    @{ Cell = 'E2'; Value = '  +0.37%  '; ForceText = $false },
    @{ Cell = 'D3'; Value = '2.664.86'; ForceText = $false },
    @{ Cell = 'E3'; Value = '  -0.05%  '; ForceText = $false },
    @{ Cell = 'E4'; Value = '  -0.05%  '; ForceText = $false },
    @{ Cell = 'D5'; Value = '598.40'; ForceText = $true },
    @{ Cell = 'E5'; Value = '  -0.32%  '; ForceText = $false },
    @{ Cell = 'D6'; Value = '157.40'; ForceText = $true },
    @{ Cell = 'E6'; Value = '  +0.30%  '; ForceText = $false },
    @{ Cell = 'D7'; Value = '0.654'; ForceText = $true },
    @{ Cell = 'E7'; Value = '  +5.13%  '; ForceText = $false },
    @{ Cell = 'D8'; Value = '0.999'; ForceText = $true },
    @{ Cell = 'E8'; Value = '  -0.06%  '; ForceText = $false },
    @{ Cell = 'D9'; Value = '0.126'; ForceText = $true },
    @{ Cell = 'E9'; Value = '  -2.50%  '; ForceText = $false },
    @{ Cell = 'D10'; Value = '0.403'; ForceText = $true },
    @{ Cell = 'E10'; Value = '  +0.22%  '; ForceText = $false },
    @{ Cell = 'D11'; Value = '5.83'; ForceText = $true },
    @{ Cell = 'E11'; Value = '  -0.22%  '; ForceText = $false },
    @{ Cell = 'E12'; Value = '  +1.57%  '; ForceText = $false },
    @{ Cell = 'D13'; Value = '28.84'; ForceText = $true },
    @{ Cell = 'E13'; Value = '  -1.36%  '; ForceText = $false },
    @{ Cell = 'D14'; Value = '0.0000193'; ForceText = $true },
    @{ Cell = 'E14'; Value = '  -2.56%  '; ForceText = $false },
    @{ Cell = 'D15'; Value = '3.142.09'; ForceText = $false },
    @{ Cell = 'E15'; Value = '  -0.11%  '; ForceText = $false },
    @{ Cell = 'D16'; Value = '65.670.09'; ForceText = $false },
    @{ Cell = 'E16'; Value = '  +0.40%  '; ForceText = $false },
    @{ Cell = 'D17'; Value = '2.631.22'; ForceText = $false },
    @{ Cell = 'E17'; Value = '  -1.53%  '; ForceText = $false },
    @{ Cell = 'D18'; Value = '12.58'; ForceText = $true },
    @{ Cell = 'E18'; Value = '  -1.38%  '; ForceText = $false },
    @{ Cell = 'D19'; Value = '4.78'; ForceText = $true },
    @{ Cell = 'E19'; Value = '  -0.14%  '; ForceText = $false },
    @{ Cell = 'D20'; Value = '348.72'; ForceText = $true },
    @{ Cell = 'D21'; Value = '7.43'; ForceText = $true },
    @{ Cell = 'E21'; Value = '  -2.07%  '; ForceText = $false },
    @{ Cell = 'D22'; Value = '0.999'; ForceText = $true },
    @{ Cell = 'E22'; Value = '  +0.00%  '; ForceText = $false },
    @{ Cell = 'D23'; Value = '69.61'; ForceText = $true },
    @{ Cell = 'E23'; Value = '  -0.02%  '; ForceText = $false },
    @{ Cell = 'D24'; Value = '1.83'; ForceText = $true },
    @{ Cell = 'E24'; Value = '  +11.21%  '; ForceText = $false },
    @{ Cell = 'D25'; Value = '0.0000112'; ForceText = $true },
    @{ Cell = 'E25'; Value = '  +0.91%  '; ForceText = $false },
    @{ Cell = 'D26'; Value = '9.57'; ForceText = $true },
    @{ Cell = 'E26'; Value = '  -0.72%  '; ForceText = $false },
    @{ Cell = 'E27'; Value = '  +2.83%  '; ForceText = $false },
    @{ Cell = 'D28'; Value = '570.30'; ForceText = $true },
    @{ Cell = 'E28'; Value = '  +7.25%  '; ForceText = $false },
    @{ Cell = 'D29'; Value = '8.12'; ForceText = $true },
    @{ Cell = 'E29'; Value = '  +0.28%  '; ForceText = $false },
    @{ Cell = 'D30'; Value = '0.163'; ForceText = $true },
    @{ Cell = 'E30'; Value = '  -2.56%  '; ForceText = $false },
    @{ Cell = 'D31'; Value = '0.999'; ForceText = $true },
    @{ Cell = 'E31'; Value = '  -0.23%  '; ForceText = $false },
    @{ Cell = 'D32'; Value = '2.14'; ForceText = $true },
    @{ Cell = 'E32'; Value = '  -0.40%  '; ForceText = $false },
    @{ Cell = 'D33'; Value = '1.82'; ForceText = $true },
    @{ Cell = 'E33'; Value = '  +3.64%  '; ForceText = $false },
    @{ Cell = 'D34'; Value = '6.69'; ForceText = $true },
    @{ Cell = 'E34'; Value = '  +4.06%  '; ForceText = $false },
    @{ Cell = 'D35'; Value = '5.45'; ForceText = $true },
    @{ Cell = 'E35'; Value = '  -0.69%  '; ForceText = $false },
    @{ Cell = 'D36'; Value = '0.422'; ForceText = $true },
    @{ Cell = 'E36'; Value = '  -0.05%  '; ForceText = $false },
    @{ Cell = 'D37'; Value = '20.58'; ForceText = $true },
    @{ Cell = 'E37'; Value = '  +0.75%  '; ForceText = $false },
    @{ Cell = 'D38'; Value = '0.999'; ForceText = $true },
    @{ Cell = 'E38'; Value = '  -0.05%  '; ForceText = $false },
    @{ Cell = 'D39'; Value = '1.94'; ForceText = $true },
    @{ Cell = 'E39'; Value = '  +0.19%  '; ForceText = $false },
    @{ Cell = 'D40'; Value = '155.40'; ForceText = $true },
    @{ Cell = 'E40'; Value = '  -1.80%  '; ForceText = $false },
    @{ Cell = 'D41'; Value = '160.50'; ForceText = $true },
    @{ Cell = 'E41'; Value = '  -2.36%  '; ForceText = $false },
    @{ Cell = 'D42'; Value = '4.08'; ForceText = $true },
    @{ Cell = 'E42'; Value = '  -1.26%  '; ForceText = $false },
    @{ Cell = 'D43'; Value = '0.0613'; ForceText = $true },
    @{ Cell = 'E43'; Value = '  +1.03%  '; ForceText = $false },
    @{ Cell = 'D44'; Value = '2.28'; ForceText = $true },
    @{ Cell = 'E44'; Value = '  -2.00%  '; ForceText = $false },
    @{ Cell = 'D45'; Value = '22.80'; ForceText = $true },
    @{ Cell = 'E45'; Value = '  -0.41%  '; ForceText = $false },
    @{ Cell = 'D46'; Value = '0.640'; ForceText = $true },
    @{ Cell = 'E46'; Value = '  -0.15%  '; ForceText = $false },
    @{ Cell = 'B47'; Value = 'VeChain'; ForceText = $false },
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; ForceText = $false },
    @{ Cell = 'D47'; Value = '0.0255'; ForceText = $true },
    @{ Cell = 'E47'; Value = '  -1.02%  '; ForceText = $false },
    @{ Cell = 'B48'; Value = 'Stellar'; ForceText = $false },
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; ForceText = $false },
    @{ Cell = 'D48'; Value = '0.102'; ForceText = $true },
    @{ Cell = 'E48'; Value = '  +1.21%  '; ForceText = $false },
    @{ Cell = 'D49'; Value = '19.81'; ForceText = $true },
    @{ Cell = 'E49'; Value = '  -1.08%  '; ForceText = $false },
    @{ Cell = 'D50'; Value = '0.0₆0243'; ForceText = $false },
    @{ Cell = 'E50'; Value = '  -2.97%  '; ForceText = $false },
    @{ Cell = 'D51'; Value = '0.804'; ForceText = $true },
    @{ Cell = 'E51'; Value = '  -1.83%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Temporarily force Text format so the numeric-looking string is not
        # reinterpreted as a Number, then drop the formatting override again so
        # the cell's style stays exactly as it was before this edit.
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.ClearFormats()
    } else {
        $rng.Value = $u.Value
    }
}

Write-Host "Applied $($updates.Count) cell updates"
